# Weekly update: insert two new price records (2023-08-16) at the top of the
# data block (rows 469/470), pushing the existing rows down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current row 469 (pushes 469..525 -> 471..527).
$ws.Rows.Item(469).EntireRow.Insert()
$ws.Rows.Item(469).EntireRow.Insert()

# New row 469 - "Primera" quality record for 2023-08-16.
$ws.Cells.Item(469, 1).Value2 = 11
$ws.Cells.Item(469, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(469, 3).Value2 = 'Bíobío'
$ws.Cells.Item(469, 4).Value2 = 45154
$ws.Cells.Item(469, 5).Value2 = 8
$ws.Cells.Item(469, 6).Value2 = 100112017
$ws.Cells.Item(469, 7).Value2 = 'Apio'
$ws.Cells.Item(469, 8).Value2 = 'Americana (o)'
$ws.Cells.Item(469, 9).Value2 = 'Primera'
$ws.Cells.Item(469, 10).Value2 = 250
$ws.Cells.Item(469, 11).Value2 = 7000
$ws.Cells.Item(469, 12).Value2 = 7500
$ws.Cells.Item(469, 13).Value2 = 7200
$ws.Cells.Item(469, 14).Value2 = '$/docena de matas'
$ws.Cells.Item(469, 15).Value2 = 'Región de Coquimbo'
$ws.Cells.Item(469, 16).Value2 = 1200
$ws.Cells.Item(469, 17).Value2 = 6
$ws.Cells.Item(469, 18).Value2 = 'Hortaliza'

# New row 470 - "Segunda" quality record for 2023-08-16.
$ws.Cells.Item(470, 1).Value2 = 11
$ws.Cells.Item(470, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(470, 3).Value2 = 'Bíobío'
$ws.Cells.Item(470, 4).Value2 = 45154
$ws.Cells.Item(470, 5).Value2 = 8
$ws.Cells.Item(470, 6).Value2 = 100112017
$ws.Cells.Item(470, 7).Value2 = 'Apio'
$ws.Cells.Item(470, 8).Value2 = 'Americana (o)'
$ws.Cells.Item(470, 9).Value2 = 'Segunda'
$ws.Cells.Item(470, 10).Value2 = 150
$ws.Cells.Item(470, 11).Value2 = 6500
$ws.Cells.Item(470, 12).Value2 = 6500
$ws.Cells.Item(470, 13).Value2 = 6500
$ws.Cells.Item(470, 14).Value2 = '$/docena de matas'
$ws.Cells.Item(470, 15).Value2 = 'Región de Coquimbo'
$ws.Cells.Item(470, 16).Value2 = 1083
$ws.Cells.Item(470, 17).Value2 = 6
$ws.Cells.Item(470, 18).Value2 = 'Hortaliza'
